$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17. This shifts the existing rows 17-35
# down to rows 18-36 (carrying their values/styles with them), matching
# the dimension growing from A1:R35 to A1:R36.
$ws.Rows.Item(17).Insert()

# The newly inserted (blank) row 17 should start out as a duplicate of
# row 16 (same market/product/region/etc.), so copy row 16 into row 17.
$ws.Range("A16:R16").Copy()
$ws.Range("A17:R17").PasteSpecial()

# Now update row 16 with this week's new price entry (weekly update per
# the commit message): new date and new min/max/avg/kg prices.
$ws.Cells.Item(16, 4).Value = 44664
$ws.Cells.Item(16, 11).Value = 1300
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = 1400
$ws.Cells.Item(16, 16).Value = 1400
